$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Move the "floors" section (previously at D15:E25) up next to the
#     "walls" section, into columns G (header) and G4:G13 (values), then
#     remove the now-empty rows 15:25. ---

# Capture the "floors" header text (D15) and the 10 item labels (E16:E25)
# before we start overwriting anything.
$floorsHeader = $ws.Cells.Item(15, 4).Value2
$floorsItems = @()
for ($i = 0; $i -lt 10; $i++) {
    $floorsItems += , $ws.Cells.Item(16 + $i, 5).Value2
}

# Capture the "walls" header text (D3) too, since it is moving from D3 to E3.
$wallsHeader = $ws.Cells.Item(3, 4).Value2

# Clear out the old "floors" block entirely (rows 15-25, columns D:E).
$ws.Range("D15:E25").ClearContents()

# --- Rebuild row 3: headers for the two side-by-side tables. ---
$ws.Cells.Item(3, 4).ClearContents()
$ws.Cells.Item(3, 5).Value = $wallsHeader
$ws.Cells.Item(3, 7).Value = $floorsHeader

# --- Renumber the D column id's for the walls table from 1-10 to 11-20. ---
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(4 + $i, 4).Value = 11 + $i
}

# --- Place the floors items alongside the walls table, in column G. ---
for ($i = 0; $i -lt 10; $i++) {
    $ws.Cells.Item(4 + $i, 7).Value = $floorsItems[$i]
}

# Row 14, column D is left as an empty, right-aligned spacer cell (matching
# the style used for the empty D3 header cell).
$ws.Cells.Item(3, 4).HorizontalAlignment = -4152
$ws.Cells.Item(14, 4).HorizontalAlignment = -4152

# Narrow spacer column F between the two tables.
$ws.Columns("F").ColumnWidth = 4.7

# Restore the selection to where the author last left the cursor.
$ws.Range("E18").Select()
